# Auto-generated Excel COM-interop script
# Applies cached-value updates (pricing data refresh) to the Famfrit_Profits workbook
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), columns H-N.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising / Quicksilver
$ws.Range("H2").Value = 599
$ws.Range("I2").Value = 619
$ws.Range("J2").Value = 499
$ws.Range("K2").Value = 619
$ws.Range("L2").Value = 499
$ws.Range("M2").Value = -506

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 2046175.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2046175.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6138526.5
$ws.Range("N17").Value = -6138862.5

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 8000
$ws.Range("I40").Value = 9333.333000000001
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 9333.333000000001
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -9158.333000000001
$ws.Range("N40").Value = -4350

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 6685.4287
$ws.Range("I62").Value = 799
$ws.Range("J62").Value = 7666.5
$ws.Range("K62").Value = 799
$ws.Range("L62").Value = 7666.5
$ws.Range("M62").Value = -175

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 6685.4287
$ws.Range("I65").Value = 799
$ws.Range("J65").Value = 7666.5
$ws.Range("K65").Value = 3995
$ws.Range("L65").Value = 38332.5
$ws.Range("M65").Value = -875

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 5014.7856
$ws.Range("I86").Value = 4201.3335
$ws.Range("J86").Value = 5624.875
$ws.Range("K86").Value = 4201.3335
$ws.Range("L86").Value = 5624.875
$ws.Range("M86").Value = -3078.3335
$ws.Range("N86").Value = -7870.875

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 5014.7856
$ws.Range("I89").Value = 4201.3335
$ws.Range("J89").Value = 5624.875
$ws.Range("K89").Value = 21006.6675
$ws.Range("L89").Value = 28124.375
$ws.Range("M89").Value = -15390.6675
$ws.Range("N89").Value = -39356.375

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 5980.5
$ws.Range("I113").Value = 5822
$ws.Range("J113").Value = 6202.4
$ws.Range("K113").Value = 5822
$ws.Range("L113").Value = 6202.4
$ws.Range("M113").Value = -2568
$ws.Range("N113").Value = -12710.4

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2400.7896
$ws.Range("I132").Value = 2616.5625
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 7849.6875
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -5319.6875

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 12206.333
$ws.Range("I137").Value = 2440.8572
$ws.Range("J137").Value = 25878
$ws.Range("K137").Value = 7322.571599999999
$ws.Range("L137").Value = 77634
$ws.Range("M137").Value = -4772.571599999999

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5704.918
$ws.Range("I138").Value = 1196.7142
$ws.Range("J138").Value = 8071.725
$ws.Range("K138").Value = 3590.1426
$ws.Range("L138").Value = 24215.175
$ws.Range("M138").Value = 1549.8574
$ws.Range("N138").Value = -34495.175

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 134.81818
$ws.Range("I5").Value = 128.3
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 128.3
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -16.30000000000001
$ws.Range("N5").Value = -424

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 52633344
$ws.Range("I61").Value = 52633344
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 52633344
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -52633132

# Row 135: Forgiveness for My Shins / Ruthenium Sabatons of Fending
$ws.Range("H135").Value = 125415
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 125415
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 125415
$ws.Range("N135").Value = -135555

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 52633344
$ws.Range("I136").Value = 52633344
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 157900032
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -157897482

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 134.81818
$ws.Range("I4").Value = 128.3
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 128.3
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -13.30000000000001
$ws.Range("N4").Value = -430

# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 4658.6665
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 4658.6665
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 4658.6665
$ws.Range("N20").Value = -5152.6665
$ws.Range("M20").ClearContents()

# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 387.375

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 2757.55
$ws.Range("I107").Value = 1493.1428
$ws.Range("J107").Value = 5707.8335
$ws.Range("K107").Value = 1493.1428
$ws.Range("L107").Value = 5707.8335
$ws.Range("M107").Value = 426.8571999999999

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 4692.077
$ws.Range("I134").Value = 4545.273
$ws.Range("J134").Value = 5499.5
$ws.Range("K134").Value = 13635.819
$ws.Range("L134").Value = 16498.5
$ws.Range("M134").Value = -11100.819

# Row 141: Awl Dreams Come True / Ra'Kaznar Awl
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 4551.1665
$ws.Range("I22").Value = 6937.933
$ws.Range("J22").Value = 573.2222
$ws.Range("K22").Value = 6937.933
$ws.Range("L22").Value = 573.2222
$ws.Range("M22").Value = -6587.933
$ws.Range("N22").Value = -1273.2222

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2048.8948
$ws.Range("I58").Value = 1808.5
$ws.Range("J58").Value = 2722
$ws.Range("K58").Value = 1808.5
$ws.Range("L58").Value = 2722
$ws.Range("M58").Value = -1605.5
$ws.Range("N58").Value = -3128

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 9930.483
$ws.Range("I99").Value = 5692.4287
$ws.Range("J99").Value = 11166.583
$ws.Range("K99").Value = 5692.4287
$ws.Range("L99").Value = 11166.583
$ws.Range("M99").Value = -4194.4287

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 9930.483
$ws.Range("I126").Value = 5692.4287
$ws.Range("J126").Value = 11166.583
$ws.Range("K126").Value = 17077.2861
$ws.Range("L126").Value = 33499.749
$ws.Range("M126").Value = -14607.2861

# Row 131: An Integral Reward / Integral Necklace of Crafting
$ws.Range("H131").Value = 21091.75
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 21091.75
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 21091.75
$ws.Range("N131").Value = -31171.75

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2048.8948
$ws.Range("I136").Value = 1808.5
$ws.Range("J136").Value = 2722
$ws.Range("K136").Value = 5425.5
$ws.Range("L136").Value = 8166
$ws.Range("M136").Value = -2875.5
$ws.Range("N136").Value = -13266

# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 106552.336
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 116121.375
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 116121.375
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -126481.375

$ws = $wb.Worksheets.Item("CUL")
# Row 18: Fisher of Men / Salt Cod
$ws.Range("H18").Value = 4163.3335
$ws.Range("I18").Value = 1245
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 3735
$ws.Range("L18").Value = 30000
$ws.Range("M18").Value = -3566

# Row 55: Pagan Pastries / Pastry Fish
$ws.Range("H55").Value = 3366.9092
$ws.Range("I55").Value = 1823
$ws.Range("J55").Value = 5219.6
$ws.Range("K55").Value = 5469
$ws.Range("L55").Value = 15658.8
$ws.Range("M55").Value = -5292
$ws.Range("N55").Value = -16012.8

# Row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 3033
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 3033
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 9099
$ws.Range("N133").Value = -19219

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 1163044.5
$ws.Range("I2").Value = 2500060.8
$ws.Range("J2").Value = 421.52173
$ws.Range("K2").Value = 2500060.8
$ws.Range("L2").Value = 421.52173
$ws.Range("M2").Value = -2499947.8

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 5119.25
$ws.Range("I80").Value = 7037.4287
$ws.Range("J80").Value = 4086.3845
$ws.Range("K80").Value = 7037.4287
$ws.Range("L80").Value = 4086.3845
$ws.Range("M80").Value = -6039.4287

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 5119.25
$ws.Range("I83").Value = 7037.4287
$ws.Range("J83").Value = 4086.3845
$ws.Range("K83").Value = 35187.14350000001
$ws.Range("L83").Value = 20431.9225
$ws.Range("M83").Value = -30195.14350000001

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 3215.926
$ws.Range("I113").Value = 2383.7334
$ws.Range("J113").Value = 4256.1665
$ws.Range("K113").Value = 2383.7334
$ws.Range("L113").Value = 4256.1665
$ws.Range("M113").Value = -213.7334000000001

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 7010.6875
$ws.Range("I132").Value = 7051.615
$ws.Range("J132").Value = 6833.3335
$ws.Range("K132").Value = 21154.845
$ws.Range("L132").Value = 20500.0005
$ws.Range("M132").Value = -18624.845
$ws.Range("N132").Value = -25560.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 2903.88
$ws.Range("I7").Value = 3078.7778
$ws.Range("J7").Value = 2454.1428
$ws.Range("K7").Value = 3078.7778
$ws.Range("L7").Value = 2454.1428
$ws.Range("M7").Value = -2966.7778

# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 1854.3334
$ws.Range("I16").Value = 1865.2
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 1865.2
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -1695.2

# Row 20: Choke Hold / Hard Leather Choker
$ws.Range("H20").Value = 12000
$ws.Range("I20").Value = 12000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -11774
$ws.Range("N20").ClearContents()

# Row 43: Subordinate Clause / Goatskin Choker
$ws.Range("H43").Value = 20749.75
$ws.Range("I43").Value = 21500
$ws.Range("J43").Value = 19999.5
$ws.Range("K43").Value = 21500
$ws.Range("L43").Value = 19999.5
$ws.Range("M43").Value = -21307
$ws.Range("N43").Value = -20385.5

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 2373.5557
$ws.Range("I68").Value = 1708.8572
$ws.Range("J68").Value = 4700
$ws.Range("K68").Value = 1708.8572
$ws.Range("L68").Value = 4700
$ws.Range("M68").Value = -959.8571999999999
$ws.Range("N68").Value = -6198

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 2373.5557
$ws.Range("I71").Value = 1708.8572
$ws.Range("J71").Value = 4700
$ws.Range("K71").Value = 8544.286
$ws.Range("L71").Value = 23500
$ws.Range("M71").Value = -4800.286
$ws.Range("N71").Value = -30988

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 2313.5
$ws.Range("I82").Value = 2223.5
$ws.Range("J82").Value = 2433.5
$ws.Range("K82").Value = 2223.5
$ws.Range("L82").Value = 2433.5
$ws.Range("M82").Value = -1862.5
$ws.Range("N82").Value = -3155.5

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 2313.5
$ws.Range("I85").Value = 2223.5
$ws.Range("J85").Value = 2433.5
$ws.Range("K85").Value = 2223.5
$ws.Range("L85").Value = 2433.5
$ws.Range("M85").Value = -975.5
$ws.Range("N85").Value = -4929.5

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 2979911.2
$ws.Range("I122").Value = 3362.9092
$ws.Range("J122").Value = 13893922
$ws.Range("K122").Value = 10088.7276
$ws.Range("L122").Value = 41681766
$ws.Range("M122").Value = -7638.7276
$ws.Range("N122").Value = -41686666

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 2903.88
$ws.Range("I126").Value = 3078.7778
$ws.Range("J126").Value = 2454.1428
$ws.Range("K126").Value = 9236.3334
$ws.Range("L126").Value = 7362.428400000001
$ws.Range("M126").Value = -6766.3334

# Row 131: For What Was Gleaned / Ophiotauroskin Wristband of Gathering
$ws.Range("H131").Value = 77268.14
$ws.Range("I131").Value = 52648
$ws.Range("J131").Value = 87116.2
$ws.Range("K131").Value = 52648
$ws.Range("L131").Value = 87116.2
$ws.Range("M131").Value = -47608
$ws.Range("N131").Value = -97196.2

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 37499.5
$ws.Range("I132").Value = 37499.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 112498.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -109968.5
